$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 728127757485.9529
    3  = 248947206859.0277
    4  = 35842103923.58646
    5  = 32789094089.75835
    6  = 24185017211.16008
    7  = 13352606143.86956
    8  = 10759940815.95023
    9  = 8966284996.765963
    10 = 8292017669.368899
    11 = 8086467521.47689
    12 = 7532202354.016511
    13 = 7190147621.414537
    14 = 6789191812.222837
    15 = 6059270581.208712
    16 = 5056212888.394469
    17 = 4774750034.886718
    18 = 4361698573.192698
    19 = 3757389910.898662
    20 = 3492021508.260048
    21 = 3330599207.781488
    22 = 3300542812.614562
    23 = 2999817564.762733
    24 = 2847168829.225609
    25 = 2750220253.895219
    26 = 2392904552.53546
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
